$wb = $excel.ActiveWorkbook

# --- Rename the "name" field to "plot_name" throughout the form ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("D3").Value = "plot_name"

$settings = $wb.Worksheets.Item("settings")
$settings.Range("B5").Value = "plot_name"

# --- Update selections on each sheet ---
$survey.Range("D24").Select() | Out-Null
$settings.Range("B5").Select() | Out-Null
